$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Rename rain model constants
$ws.Range("A19").Value = "fi_lidar_rain_reflectivity"
$ws.Range("A20").Value = "fi_lidar_rain_intensity"

# Update selection to reflect where the cursor ended up
$ws.Range("A20").Select()
